$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table has a single column; each row holds one measurement value in a
# single run. A handful of rows get their value swapped for a new one, two
# new rows' worth of values are folded in (row count stays the same because
# three trailing rows collapse from "tab separated list of 10 values in one
# run" down to "a single value"), and the trailing three rows inherit the
# values that used to live in rows 1-3.

function Set-CellText($table, $row, $text) {
    $table.Cell($row, 1).Range.Text = $text
}

Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "322"
Set-CellText $t 5  "0.00002"
Set-CellText $t 6  "0.00016"
# row 7 ("0.00005") is unchanged
Set-CellText $t 8  "0.00003"
Set-CellText $t 9  "0.00004"
# row 10 ("0.00005") is unchanged
Set-CellText $t 11 "0.00010"
Set-CellText $t 12 "0.01560"

# Rows 44-46 currently hold 10 tab-separated values inside a single run;
# replace the whole cell contents with the single surviving value.
Set-CellText $t 44 "99.98"
Set-CellText $t 45 "0.02"
Set-CellText $t 46 "72"

Write-Host "Done. Table now has" $t.Rows.Count "rows."
